$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$s.Shapes.Item("Picture 9").Delete()
